# TP3: add two GitHub project-link hyperlinks after the two Heading2
# "exercise" paragraphs, each followed by the original blank paragraph
# that was already there (matching the authoring diff).

$d = $word.ActiveDocument
$url = "https://github.com/faculdade-infnet/V-2-Microsservicos-e-Spring-Cloud/tree/main/TP3/TP3-projeto"

function Get-ParagraphIndexAt($doc, $pos) {
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Start -eq $pos) {
            return $i
        }
    }
    return -1
}

function Add-LinkParagraphAfterHeading($doc, $headingText, $linkUrl) {
    # Locate the heading paragraph by its exact visible text.
    $found = $doc.Content
    $ok = $found.Find.Execute($headingText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $null = $found.Expand(4)   # wdParagraph -> grow to include the paragraph mark

    $insertPos = $found.End

    # Insert a brand-new (Normal-styled) empty paragraph right after the
    # heading; this pushes the pre-existing blank paragraph one slot later.
    $ins = $doc.Range($insertPos, $insertPos)
    $ins.InsertParagraphAfter()

    $newIdx = Get-ParagraphIndexAt $doc $insertPos
    $newPara = $doc.Paragraphs.Item($newIdx)
    $newPara.Range.InsertAfter($linkUrl)

    $newPara2 = $doc.Paragraphs.Item($newIdx)
    $textRange = $doc.Range($newPara2.Range.Start, $newPara2.Range.End - 1)
    $null = $doc.Hyperlinks.Add($textRange, $linkUrl, "", "", $linkUrl)

    return $newIdx
}

# --- Link #1: right after "Implemente uma aplicação ..." ---
$heading1 = "Implemente uma aplicação seguindo os princípios reativos, adicione persistência usando o Spring Data JDBC, use o cliente WebClient para fazer as requisições HTTP e faça os testes usando Testcontainers."
$idx1 = Add-LinkParagraphAfterHeading $d $heading1 $url

# --- Link #2: right after "Para a aplicação solicitada ..." ---
# An extra blank paragraph is also inserted after this second link (per the
# target diff), on top of the pre-existing blank paragraph that follows it.
$heading2 = "Para a aplicação solicitada, vocês devem usar o conceito de microsserviços e seus principais recursos, como Docker e Kubernetes."
$idx2 = Add-LinkParagraphAfterHeading $d $heading2 $url

$linkPara2 = $d.Paragraphs.Item($idx2)
$afterLink2 = $linkPara2.Range.End
$insBlank = $d.Range($afterLink2, $afterLink2)
$insBlank.InsertParagraphAfter()

Write-Output ("Done. Hyperlinks.Count=" + $d.Hyperlinks.Count + " InlineShapes.Count=" + $d.InlineShapes.Count)
